$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "Despesa"
$ws.Range("B12").Value = "Retirada de Capital"
$ws.Range("C12").Value = 100
$ws.Range("D12").Value = "19/02/2025"
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = "Luiz"
